$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.046.58"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.589.00"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'578.37"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").Value = "'190.92"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").Value = "3.589.58"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "'0.665"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "'55.83"
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("D13").Value = "'0.0000307"
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "4.162.21"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "3.585.28"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "69.963.79"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "'474.72"
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").Value = "'19.22"
$ws.Range("E23").Value = "  +11.02%  "
$ws.Range("E24").Value = "  -6.50%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'4.37"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'95.27"
$ws.Range("E26").Value = "  +4.71%  "
$ws.Range("D27").Value = "'2.99"
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("D28").Value = "'11.02"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").Value = "'32.26"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "'7.65"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "'66.49"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "'589.00"
$ws.Range("E35").Value = "  -6.00%  "
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "0.0₃0802"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  +17.31%  "
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("D42").Value = "'3.45"
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("D43").Value = "3.220.17"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("D45").Value = "'3.06"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "'0.0445"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "'9.45"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  -4.95%  "
